$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the _GoBack bookmark from its original spot (the empty
#    paragraph right after "N SIRET :" in the header block).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Collapse the closing-clause paragraphs:
#      [spacing 263, empty]
#      "En votre aimable r\u00e8glement,"
#      [spacing 30, empty]
#      "Cordialement,"
#      [spacing 249, empty]
#    down to a single empty paragraph (keeping the first one, whose
#    spacing becomes 249) right before "Conditions de paiement".
# ------------------------------------------------------------------
$pCordialement = $null
$pAimable = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Cordialement*") { $pCordialement = $p }
    if ($t -like "*aimable*") { $pAimable = $p }
}

$pKeep = $pAimable.Previous()        # the "spacing=263" empty paragraph to keep
$pLast = $pCordialement.Next()       # the "spacing=249" empty paragraph to drop

# Remove everything from just after $pKeep's paragraph mark through the
# end of $pLast's paragraph mark - i.e. the 4 in-between paragraphs.
$cutRange = $d.Range($pKeep.Range.End, $pLast.Range.End)
$cutRange.Delete()

# The kept paragraph now carries the final 249-twip (12.45pt) spacing.
$pKeep.Format.LineSpacing = 12.45

# ------------------------------------------------------------------
# 3) Re-create the _GoBack bookmark inside that now-merged, still
#    text-less paragraph. A genuinely zero-length range anchored in
#    an empty paragraph does not seat correctly, so stash a throwaway
#    character, bookmark across it, then delete the character again -
#    the bookmark collapses to the right spot and survives.
# ------------------------------------------------------------------
$pKeep.Range.InsertBefore("x")
$bmRange = $d.Range($pKeep.Range.Start, $pKeep.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$tmpRange = $d.Range($pKeep.Range.Start, $pKeep.Range.Start + 1)
$tmpRange.Delete()
